# "Generate Report for Handoff"
#
# The localization report previously showed 348fa1fa-...md as the first
# (row 2) entry and fa6cc81e-...md as the second (row 3) entry on every
# sheet, both still sitting at status "Handed back: in sync with en-US".
#
# After re-generating the report for handoff:
#   - fa6cc81e-...md now sorts first (row 2) and keeps the handed-back
#     status.
#   - 348fa1fa-...md now sorts second (row 3), has moved on to
#     "Ready for handoff" with fresh handoff timestamps, and (on the
#     per-locale sheets) carries a stale-handback warning in the
#     "Error Detail" column, whose column got widened to fit the text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
$ws.Range("A3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"

$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-01 16:57:46"

# Hyperlink display text (B2/B3) swaps along with the row contents.
$idx = 0
foreach ($h in $ws.Hyperlinks) {
  $idx = $idx + 1
  if ($idx -eq 1) {
    $h.TextToDisplay = "e2e\fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
  }
  if ($idx -eq 2) {
    $h.TextToDisplay = "e2e\348fa1fa-b11a-4d9d-8129-e29380aef063.md"
  }
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
$ws.Range("G2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.5c42f61c24227defa8105b0dc1201e6444cd2316.zh-cn.xlf"
$ws.Range("I2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
$ws.Range("J2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.5c42f61c24227defa8105b0dc1201e6444cd2316.zh-cn.xlf"

$ws.Range("A3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.cf6e5663733644c5bd83e1664d522f891fb3e96d.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-01 16:57:42"
$ws.Range("I3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
$ws.Range("J3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.cf6e5663733644c5bd83e1664d522f891fb3e96d.zh-cn.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62f83ce84f091623428b2b0d7894e3fa12352346/e2e/348fa1fa-b11a-4d9d-8129-e29380aef063.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d8d7e86d0b114f131cebc2418fb16c0e520af99/e2e/348fa1fa-b11a-4d9d-8129-e29380aef063.md."

# Column P ("Error Detail") widened to fit the new message.
$ws.Columns.Item(16).ColumnWidth = 39.17

$idx = 0
foreach ($h in $ws.Hyperlinks) {
  $idx = $idx + 1
  if ($idx -eq 1) {
    $h.TextToDisplay = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
  }
  if ($idx -eq 2) {
    $h.TextToDisplay = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
  }
  if ($idx -eq 3) {
    $h.TextToDisplay = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
  }
  if ($idx -eq 4) {
    $h.TextToDisplay = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
  }
}

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
$ws.Range("G2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.5c42f61c24227defa8105b0dc1201e6444cd2316.de-de.xlf"
$ws.Range("I2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
$ws.Range("J2").Value = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.5c42f61c24227defa8105b0dc1201e6444cd2316.de-de.xlf"

$ws.Range("A3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.cf6e5663733644c5bd83e1664d522f891fb3e96d.de-de.xlf"
$ws.Range("H3").Value = "2016-09-01 16:57:46"
$ws.Range("I3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
$ws.Range("J3").Value = "348fa1fa-b11a-4d9d-8129-e29380aef063.cf6e5663733644c5bd83e1664d522f891fb3e96d.de-de.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62f83ce84f091623428b2b0d7894e3fa12352346/e2e/348fa1fa-b11a-4d9d-8129-e29380aef063.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d8d7e86d0b114f131cebc2418fb16c0e520af99/e2e/348fa1fa-b11a-4d9d-8129-e29380aef063.md."

# Column P ("Error Detail") widened to fit the new message.
$ws.Columns.Item(16).ColumnWidth = 39.17

$idx = 0
foreach ($h in $ws.Hyperlinks) {
  $idx = $idx + 1
  if ($idx -eq 1) {
    $h.TextToDisplay = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
  }
  if ($idx -eq 2) {
    $h.TextToDisplay = "fa6cc81e-9dee-4c3a-a716-11e11a9f3b43.md"
  }
  if ($idx -eq 3) {
    $h.TextToDisplay = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
  }
  if ($idx -eq 4) {
    $h.TextToDisplay = "348fa1fa-b11a-4d9d-8129-e29380aef063.md"
  }
}
